$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "James" -> "James Kenniff" (name correction)
$ws.Range("A2").Value = "James Kenniff"

# 4 new data points (rows 3-6)
$ws.Range("A3").Value = "Kiki Apple"
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = "Max Gallo"
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "Used the hard punch to clear his picture without instruction"

$ws.Range("A5").Value = "Andy Yang"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("A6").Value = "Andy Liu"
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0

$ws.Range("A7").Select()
